$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.366.02"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.523.12"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.67"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.39"
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.05"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.65"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "2.909.81"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "2.517.74"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.860"
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "42.490.93"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.59"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.79"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.65"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.19"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.14"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.92"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.26"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.12"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.30"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.80"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0782"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.118"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.77"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.34"
$ws.Range("E41").Value = "  +11.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").Value = "2.025.12"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.06"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.83"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").Value = "2.766.67"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.85"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -1.79%  "
